$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B27").Value = "ThingDef"
$ws.Range("E27").Value = "Install {0_label} to become mechanitor"
$ws.Range("F27").Value = "메카나이터가 되기 위해 {0_label} 설치"
$ws.Range("F27").Font.Name = "맑은 고딕"
$ws.Range("F27").Font.Size = 11
$ws.Range("F27").Font.Family = 3
$ws.Range("C27").Value = "Mechlink.comps.CompUsableImplant_NoMechanitor.useLabel"
$ws.Range("A27").Value = "ThingDef+Mechlink.comps.CompUsableImplant_NoMechanitor.useLabel"

$ws.Range("E28").Select()
